$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.515.27"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "3.690.08"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'677.89"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "'161.54"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'7.14"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("D11").Value = "'0.439"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "'0.0000233"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "4.313.50"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Value = "'32.48"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "3.682.65"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "69.452.87"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "'0.116"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").Value = "'16.03"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'6.47"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "'470.74"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").Value = "3.837.38"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "'10.87"
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("D28").Value = "'9.13"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").Value = "'2.70"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "'6.59"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("D34").Value = "'27.00"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "3.680.46"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").Value = "'8.47"
$ws.Range("E37").Value = "  +3.04%  "
$ws.Range("D38").Value = "'6.23"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.28"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'0.0901"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").Value = "'168.64"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").Value = "'0.943"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "'46.70"
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "'0.000278"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "'27.63"
$ws.Range("E49").Value = "  -3.91%  "
$ws.Range("D50").Value = "'1.09"
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("D51").Value = "'7.88"
$ws.Range("E51").Value = "  +0.74%  "
